$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 90
$ws.Range("A90").Value = 14601531
$ws.Range("B90").NumberFormat = "@"
$ws.Range("B90").Value = '2025-09-12'
$ws.Range("B90").ClearFormats()
$ws.Range("C90").Value = 'Janice Tjen'
$ws.Range("D90").Value = 'Alexandra Eala'
$ws.Range("E90").Value = 'Gana Alexandra Eala'
$ws.Range("F90").Value = 1.8
$ws.Range("G90").Value = ""
$ws.Range("G90").Font.Bold = $false
$ws.Range("G90").ClearFormats()
$ws.Range("H90").Value = ""
$ws.Range("H90").Font.Bold = $false
$ws.Range("H90").ClearFormats()

# Row 91
$ws.Range("A91").Value = 14601543
$ws.Range("B91").NumberFormat = "@"
$ws.Range("B91").Value = '2025-09-12'
$ws.Range("B91").ClearFormats()
$ws.Range("C91").Value = 'Francesca Jones'
$ws.Range("D91").Value = 'Solana Sierra'
$ws.Range("E91").Value = 'Gana Solana Sierra'
$ws.Range("F91").Value = 2.3
$ws.Range("G91").Value = ""
$ws.Range("G91").Font.Bold = $false
$ws.Range("G91").ClearFormats()
$ws.Range("H91").Value = ""
$ws.Range("H91").Font.Bold = $false
$ws.Range("H91").ClearFormats()

# Row 92
$ws.Range("A92").Value = 14601414
$ws.Range("B92").NumberFormat = "@"
$ws.Range("B92").Value = '2025-09-12'
$ws.Range("B92").ClearFormats()
$ws.Range("C92").Value = 'Patrick Zahraj'
$ws.Range("D92").Value = 'Cannon Kingsley'
$ws.Range("E92").Value = 'Gana Cannon Kingsley'
$ws.Range("F92").Value = 2.2
$ws.Range("G92").Value = ""
$ws.Range("G92").Font.Bold = $false
$ws.Range("G92").ClearFormats()
$ws.Range("H92").Value = ""
$ws.Range("H92").Font.Bold = $false
$ws.Range("H92").ClearFormats()

# Row 93
$ws.Range("A93").Value = 14601394
$ws.Range("B93").NumberFormat = "@"
$ws.Range("B93").Value = '2025-09-12'
$ws.Range("B93").ClearFormats()
$ws.Range("C93").Value = 'Stefano Napolitano'
$ws.Range("D93").Value = 'Pol Martin Tiffon'
$ws.Range("E93").Value = 'Gana Stefano Napolitano'
$ws.Range("F93").Value = 1.57
$ws.Range("G93").Value = ""
$ws.Range("G93").Font.Bold = $false
$ws.Range("G93").ClearFormats()
$ws.Range("H93").Value = ""
$ws.Range("H93").Font.Bold = $false
$ws.Range("H93").ClearFormats()

# Row 94
$ws.Range("A94").Value = 14604885
$ws.Range("B94").NumberFormat = "@"
$ws.Range("B94").Value = '2025-09-12'
$ws.Range("B94").ClearFormats()
$ws.Range("C94").Value = 'Filip Cristian Jianu'
$ws.Range("D94").Value = 'Stefan Palosi'
$ws.Range("E94").Value = 'Gana Stefan Palosi'
$ws.Range("F94").Value = 3.4
$ws.Range("G94").Value = ""
$ws.Range("G94").Font.Bold = $false
$ws.Range("G94").ClearFormats()
$ws.Range("H94").Value = ""
$ws.Range("H94").Font.Bold = $false
$ws.Range("H94").ClearFormats()

# Row 95
$ws.Range("A95").Value = 14601429
$ws.Range("B95").NumberFormat = "@"
$ws.Range("B95").Value = '2025-09-12'
$ws.Range("B95").ClearFormats()
$ws.Range("C95").Value = 'Trevor Svajda'
$ws.Range("D95").Value = 'Murphy Cassone'
$ws.Range("E95").Value = 'Gana Trevor Svajda'
$ws.Range("F95").Value = 2
$ws.Range("G95").Value = ""
$ws.Range("G95").Font.Bold = $false
$ws.Range("G95").ClearFormats()
$ws.Range("H95").Value = ""
$ws.Range("H95").Font.Bold = $false
$ws.Range("H95").ClearFormats()

# Row 96
$ws.Range("A96").Value = 14601425
$ws.Range("B96").NumberFormat = "@"
$ws.Range("B96").Value = '2025-09-12'
$ws.Range("B96").ClearFormats()
$ws.Range("C96").Value = 'Jack Pinnington Jones'
$ws.Range("D96").Value = 'Aidan Kim'
$ws.Range("E96").Value = 'Gana Aidan Kim'
$ws.Range("F96").Value = 3.25
$ws.Range("G96").Value = ""
$ws.Range("G96").Font.Bold = $false
$ws.Range("G96").ClearFormats()
$ws.Range("H96").Value = ""
$ws.Range("H96").Font.Bold = $false
$ws.Range("H96").ClearFormats()

# Row 97
$ws.Range("A97").Value = 14668918
$ws.Range("B97").NumberFormat = "@"
$ws.Range("B97").Value = '2025-09-12'
$ws.Range("B97").ClearFormats()
$ws.Range("C97").Value = 'Karim Ibrahim'
$ws.Range("D97").Value = 'Liam Broady'
$ws.Range("E97").Value = 'Gana Karim Ibrahim'
$ws.Range("F97").Value = 8
$ws.Range("G97").Value = ""
$ws.Range("G97").Font.Bold = $false
$ws.Range("G97").ClearFormats()
$ws.Range("H97").Value = ""
$ws.Range("H97").Font.Bold = $false
$ws.Range("H97").ClearFormats()

# Row 98
$ws.Range("A98").Value = 14668899
$ws.Range("B98").NumberFormat = "@"
$ws.Range("B98").Value = '2025-09-12'
$ws.Range("B98").ClearFormats()
$ws.Range("C98").Value = 'Jeremy Gschwendtner'
$ws.Range("D98").Value = 'George Lazarov'
$ws.Range("E98").Value = 'Gana George Lazarov'
$ws.Range("F98").Value = 1.83
$ws.Range("G98").Value = ""
$ws.Range("G98").Font.Bold = $false
$ws.Range("G98").ClearFormats()
$ws.Range("H98").Value = ""
$ws.Range("H98").Font.Bold = $false
$ws.Range("H98").ClearFormats()

# Row 99
$ws.Range("A99").Value = 14667149
$ws.Range("B99").NumberFormat = "@"
$ws.Range("B99").Value = '2025-09-12'
$ws.Range("B99").ClearFormats()
$ws.Range("C99").Value = 'Giles Hussey'
$ws.Range("D99").Value = 'João Loureiro'
$ws.Range("E99").Value = 'Gana João Loureiro'
$ws.Range("F99").Value = 6
$ws.Range("G99").Value = ""
$ws.Range("G99").Font.Bold = $false
$ws.Range("G99").ClearFormats()
$ws.Range("H99").Value = ""
$ws.Range("H99").Font.Bold = $false
$ws.Range("H99").ClearFormats()

# Row 100
$ws.Range("A100").Value = 14666286
$ws.Range("B100").NumberFormat = "@"
$ws.Range("B100").Value = '2025-09-12'
$ws.Range("B100").ClearFormats()
$ws.Range("C100").Value = 'Manas Dhamne'
$ws.Range("D100").Value = 'Luca Wiedenmann'
$ws.Range("E100").Value = 'Gana Luca Wiedenmann'
$ws.Range("F100").Value = 2.75
$ws.Range("G100").Value = ""
$ws.Range("G100").Font.Bold = $false
$ws.Range("G100").ClearFormats()
$ws.Range("H100").Value = ""
$ws.Range("H100").Font.Bold = $false
$ws.Range("H100").ClearFormats()

# Row 101
$ws.Range("A101").Value = 14667181
$ws.Range("B101").NumberFormat = "@"
$ws.Range("B101").Value = '2025-09-12'
$ws.Range("B101").ClearFormats()
$ws.Range("C101").Value = 'Hamish Stewart'
$ws.Range("D101").Value = 'Pierre Delage'
$ws.Range("E101").Value = 'Gana Pierre Delage'
$ws.Range("F101").Value = 3.25
$ws.Range("G101").Value = ""
$ws.Range("G101").Font.Bold = $false
$ws.Range("G101").ClearFormats()
$ws.Range("H101").Value = ""
$ws.Range("H101").Font.Bold = $false
$ws.Range("H101").ClearFormats()

# Row 102
$ws.Range("A102").Value = 14668434
$ws.Range("B102").NumberFormat = "@"
$ws.Range("B102").Value = '2025-09-12'
$ws.Range("B102").ClearFormats()
$ws.Range("C102").Value = 'Martin  Sabas'
$ws.Range("D102").Value = 'Harry Wendelken'
$ws.Range("E102").Value = 'Gana Martin  Sabas'
$ws.Range("F102").Value = 4.5
$ws.Range("G102").Value = ""
$ws.Range("G102").Font.Bold = $false
$ws.Range("G102").ClearFormats()
$ws.Range("H102").Value = ""
$ws.Range("H102").Font.Bold = $false
$ws.Range("H102").ClearFormats()

# Row 103
$ws.Range("A103").Value = 14667185
$ws.Range("B103").NumberFormat = "@"
$ws.Range("B103").Value = '2025-09-12'
$ws.Range("B103").ClearFormats()
$ws.Range("C103").Value = 'James Story'
$ws.Range("D103").Value = 'Mats Rosenkranz'
$ws.Range("E103").Value = 'Gana James Story'
$ws.Range("F103").Value = 2.25
$ws.Range("G103").Value = ""
$ws.Range("G103").Font.Bold = $false
$ws.Range("G103").ClearFormats()
$ws.Range("H103").Value = ""
$ws.Range("H103").Font.Bold = $false
$ws.Range("H103").ClearFormats()

# Row 104
$ws.Range("A104").Value = 14667168
$ws.Range("B104").NumberFormat = "@"
$ws.Range("B104").Value = '2025-09-12'
$ws.Range("B104").ClearFormats()
$ws.Range("C104").Value = 'Ryan Peniston'
$ws.Range("D104").Value = 'Marcello Serafini'
$ws.Range("E104").Value = 'Gana Marcello Serafini'
$ws.Range("F104").Value = 5
$ws.Range("G104").Value = ""
$ws.Range("G104").Font.Bold = $false
$ws.Range("G104").ClearFormats()
$ws.Range("H104").Value = ""
$ws.Range("H104").Font.Bold = $false
$ws.Range("H104").ClearFormats()

